# Adds a new person ("Federico Carusso") as row 5 of the people table on
# Sheet1, wires up the mailto: hyperlink for his e-mail cell (matching the
# style used by the existing hyperlink cells above it), and leaves the
# selection where the user ended up after entering the data (cell I6, just
# past the bottom-right of the newly-entered row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: new record ------------------------------------------------
$ws.Range("A5").Value = "Federico Carusso"
$ws.Range("B5").Value = "carusso@test.test"
$ws.Range("C5").Value = "male"
$ws.Range("D5").Value = "Italian"
$ws.Range("E5").Value = 5555555

# The birthdate column stores plain text dates (e.g. "1986.11.25") rather
# than real date serials elsewhere in the sheet, so force this cell to text
# before typing the value to stop it being auto-recognised as a date, then
# drop back to the default (unstyled) format once the text is committed.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "1966.09.15"
$ws.Range("F5").ClearFormats()

$ws.Range("G5").Value = "Italia"
$ws.Range("H5").Value = "Project Manager"
$ws.Range("I5").Value = "Sicilia"
$ws.Range("J5").Value = "Italy"
$ws.Range("K5").Value = "Pasta 8888"
$ws.Range("L5").Value = 22222
$ws.Range("M5").Value = "75 inches"

# --- Hyperlink on the e-mail cell, matching B2:B4 ----------------------
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:carusso@test.test")
$ws.Range("B5").Style = "Hyperlink"

# --- Final selection left on the sheet ---------------------------------
$ws.Range("I6").Select()
